$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data is written column-by-column (matches how the shared-string table was
# populated in the authored workbook) so new shared strings land in the
# same order as the source edit.

$textFmt = $ws.Range("A2").NumberFormat

# --- Column A: Ma so (IDs, text-formatted column) ---
$ws.Range("A2").Value = "0027"
$ws.Range("A3").Value = "0028"
$ws.Range("A4").Value = "0029"
$ws.Range("A5").Value = "0030"
$ws.Range("A6").Value = "0031"
$ws.Range("A7").NumberFormat = $textFmt
$ws.Range("A7").Value = "0032"
$ws.Range("A8").NumberFormat = $textFmt
$ws.Range("A8").Value = "0033"
$ws.Range("A9").NumberFormat = $textFmt
$ws.Range("A9").Value = "0034"
$ws.Range("A10").NumberFormat = $textFmt
$ws.Range("A10").Value = "0035"

# --- Column B: Dia chi ---
$ws.Range("B6").Value = "đồng nai"
$ws.Range("B7").Value = "đồng nai"
$ws.Range("B8").Value = "đồng nai"
$ws.Range("B9").Value = "đồng nai"
$ws.Range("B10").Value = "đồng nai"

# --- Column C: Gioi tinh ---
$ws.Range("C6").Value = "Nam"
$ws.Range("C7").Value = "Nam"
$ws.Range("C8").Value = "Nam"
$ws.Range("C9").Value = "Nam"
$ws.Range("C10").Value = "Nữ"

# --- Column D: Ho ten ---
$ws.Range("D6").Value = "Văn E"
$ws.Range("D7").Value = "Văn F"
$ws.Range("D8").Value = "Văn G"
$ws.Range("D9").Value = "Văn H"
$ws.Range("D10").Value = "Văn Q"

# --- Column E: Ngay sinh (text-formatted column) ---
$ws.Range("E6").Value = "2000-11-7"
$ws.Range("E7").NumberFormat = $textFmt
$ws.Range("E7").Value = "2000-11-8"
$ws.Range("E8").NumberFormat = $textFmt
$ws.Range("E8").Value = "2000-11-9"
$ws.Range("E9").NumberFormat = $textFmt
$ws.Range("E9").Value = "2000-11-10"
$ws.Range("E10").NumberFormat = $textFmt
$ws.Range("E10").Value = "2000-11-11"

# --- Column F: So dien thoai (mixed text/number, matching source pattern) ---
$ws.Range("F3").Value = "909526212"
# F4 already carries the text-format style from the template; clear it first
# so the numeric literal isn't coerced to text, then restore the look.
$ws.Range("F4").ClearFormats()
$ws.Range("F4").Value = 909526212
$ws.Range("F4").NumberFormat = $textFmt
$ws.Range("F5").Value = "909526213"
$ws.Range("F6").ClearFormats()
$ws.Range("F6").Value = 909526213
$ws.Range("F6").NumberFormat = $textFmt
$ws.Range("F7").NumberFormat = $textFmt
$ws.Range("F7").Value = "909526214"
$ws.Range("F8").Value = 909526214
$ws.Range("F8").NumberFormat = $textFmt
$ws.Range("F9").NumberFormat = $textFmt
$ws.Range("F9").Value = "909526215"
$ws.Range("F10").Value = 909526215
$ws.Range("F10").NumberFormat = $textFmt

# --- Column G: Khoa hoc ---
$ws.Range("G2").Value = "2018-2022"
$ws.Range("G3").Value = "2018-2022"
$ws.Range("G4").Value = "2018-2022"
$ws.Range("G5").Value = "2018-2022"
$ws.Range("G6").Value = "2018-2022"
$ws.Range("G7").Value = "2018-2022"
$ws.Range("G8").Value = "2018-2022"
$ws.Range("G9").Value = "2018-2022"
$ws.Range("G10").Value = "2018-2022"

# --- Sheet-level cosmetics ---
# Column E width changed (narrower bestFit text replaced by an explicit custom width)
$ws.Columns("E").ColumnWidth = 12.67

# Active selection moves to G12 (below the last used row)
$ws.Range("G12").Select()
